$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price/volume columns to text format so numeric-looking strings
# (e.g. "1.00", "0.999") are preserved exactly as text, not coerced to numbers.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "63.025.27"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").Value = "3.164.42"
$ws.Range("E3").Value = "  -4.63%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "591.35"
$ws.Range("E5").Value = "  -2.17%  "
$ws.Range("D6").Value = "134.25"
$ws.Range("E6").Value = "  -5.88%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.159.88"
$ws.Range("E8").Value = "  -4.71%  "
$ws.Range("D9").Value = "0.517"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("E10").Value = "  -5.61%  "
$ws.Range("D11").Value = "5.24"
$ws.Range("E11").Value = "  -5.50%  "
$ws.Range("E12").Value = "  -3.42%  "
$ws.Range("E13").Value = "  -4.27%  "
$ws.Range("D14").Value = "34.70"
$ws.Range("E14").Value = "  -0.61%  "
$ws.Range("D15").Value = "3.686.22"
$ws.Range("E15").Value = "  -4.61%  "
$ws.Range("E16").Value = "  -1.80%  "
$ws.Range("D17").Value = "3.166.40"
$ws.Range("E17").Value = "  -4.69%  "
$ws.Range("D18").Value = "62.968.13"
$ws.Range("E18").Value = "  -1.36%  "
$ws.Range("D19").Value = "6.58"
$ws.Range("E19").Value = "  -4.51%  "
$ws.Range("D20").Value = "462.13"
$ws.Range("E20").Value = "  -3.91%  "
$ws.Range("D21").Value = "13.86"
$ws.Range("E21").Value = "  -1.82%  "
$ws.Range("D22").Value = "0.699"
$ws.Range("E22").Value = "  -5.15%  "
$ws.Range("D23").Value = "7.62"
$ws.Range("E23").Value = "  -5.49%  "
$ws.Range("D24").Value = "13.38"
$ws.Range("E24").Value = "  -2.61%  "
$ws.Range("D25").Value = "83.49"
$ws.Range("E25").Value = "  -1.51%  "
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  -0.13%  "
$ws.Range("E27").Value = "  +0.14%  "
$ws.Range("E28").Value = "  -3.91%  "
$ws.Range("E29").Value = "  -6.12%  "
$ws.Range("D30").Value = "6.73"
$ws.Range("E30").Value = "  -6.67%  "
$ws.Range("D31").Value = "2.03"
$ws.Range("E31").Value = "  -6.24%  "
$ws.Range("D32").Value = "27.17"
$ws.Range("E32").Value = "  -6.32%  "
$ws.Range("E33").Value = "  -3.26%  "
$ws.Range("D34").Value = "2.36"
$ws.Range("E34").Value = "  -6.65%  "
$ws.Range("D35").Value = "1.03"
$ws.Range("E35").Value = "  -6.84%  "
$ws.Range("D36").Value = "5.83"
$ws.Range("E36").Value = "  -4.42%  "
$ws.Range("D37").Value = "51.26"
$ws.Range("E37").Value = "  -1.95%  "
$ws.Range("D38").Value = "0.0₃0705"
$ws.Range("E38").Value = "  -5.21%  "
$ws.Range("E39").Value = "  -3.08%  "
$ws.Range("D40").Value = "404.97"
$ws.Range("E40").Value = "  -7.10%  "
$ws.Range("D41").Value = "8.14"
$ws.Range("E41").Value = "  -2.42%  "
$ws.Range("E42").Value = "  -6.70%  "
$ws.Range("D43").Value = "2.60"
$ws.Range("E43").Value = "  -6.20%  "
$ws.Range("D44").Value = "2.803.55"
$ws.Range("E44").Value = "  -9.61%  "
$ws.Range("E45").Value = "  -5.46%  "
$ws.Range("E46").Value = "  +0.02%  "
$ws.Range("E47").Value = "  -6.02%  "
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").Value = "123.68"
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").Value = "25.31"
$ws.Range("E49").Value = "  -4.03%  "
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("D51").Value = "34.04"
$ws.Range("E51").Value = "  -9.55%  "

# Restore default (Normal) style/number format so the cells do not retain
# the temporary text-format styling.
$dataRange.Style = "Normal"
